$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from the existing header cell (H1) onto the two new
# header cells so they share the same style index instead of minting a
# new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells I2 and J2 (plain numeric values, no special style)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2
